$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'62.247.60"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.44%  '
$ws.Range('D3').Value = "'3.426.27"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.55%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').Value = "'407.09"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').Value = "'132.17"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.97%  '
$ws.Range('E7').Value = '  -1.88%  '
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('E9').Value = '  +3.61%  '
$ws.Range('D10').Value = "'0.137"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +10.27%  '
$ws.Range('D11').Value = "'41.86"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.14%  '
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('E13').Value = '  +1.85%  '
$ws.Range('D14').Value = "'8.40"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.85%  '
$ws.Range('D15').Value = "'3.467.05"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').Value = "'11.61"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.30%  '
$ws.Range('D17').Value = "'62.159.98"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('E18').Value = '  +0.62%  '
$ws.Range('D19').Value = "'0.0000149"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +12.07%  '
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('D21').Value = "'84.35"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.39%  '
$ws.Range('D22').Value = "'312.05"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.23%  '
$ws.Range('D23').Value = "'12.76"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.15%  '
$ws.Range('D24').Value = "'3.17"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('D25').Value = "'4.76"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.86%  '
$ws.Range('D26').Value = "'29.71"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('D27').Value = "'8.15"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.27%  '
$ws.Range('E28').Value = '  +5.38%  '
$ws.Range('D29').Value = "'2.79"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.65%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = "'44.11"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.46%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').Value = "'0.172"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.17%  '
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('D33').Value = "'11.34"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.87%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').Value = "'0.0487"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').Value = "'51.63"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('D37').Value = "'0.999"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.24%  '
$ws.Range('E38').Value = '  +1.80%  '
$ws.Range('E39').Value = '  -2.29%  '
$ws.Range('D40').Value = "'0.315"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.01%  '
$ws.Range('D41').Value = "'141.06"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.79%  '
$ws.Range('E42').Value = '  +0.35%  '
$ws.Range('D43').Value = "'1.98"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('E44').Value = '  +0.66%  '
$ws.Range('D45').Value = "'16.77"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.41%  '
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').Value = "'21.37"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.04%  '
$ws.Range('D48').Value = "'2.100.90"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('D49').Value = "'2.31"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.65%  '
$ws.Range('E50').Value = '  +2.85%  '
$ws.Range('D51').Value = "'1.71"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +18.16%  '
